# Adapt column header formatting to respective input file names, add a
# worksheet Table over the data, and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename the "_old" / "_new" header-name suffixes to the
#        format-version-specific suffixes "_FV2410" / "_FV2504". -----------
$oldToNew = @{
    "Segmentname_old"         = "Segmentname_FV2410"
    "Segmentgruppe_old"       = "Segmentgruppe_FV2410"
    "Segment_old"             = "Segment_FV2410"
    "Datenelement_old"        = "Datenelement_FV2410"
    "Segment ID_old"          = "Segment ID_FV2410"
    "Code_old"                = "Code_FV2410"
    "Qualifier_old"           = "Qualifier_FV2410"
    "Beschreibung_old"        = "Beschreibung_FV2410"
    "Bedingungsausdruck_old"  = "Bedingungsausdruck_FV2410"
    "Bedingung_old"           = "Bedingung_FV2410"
    "Segmentname_new"         = "Segmentname_FV2504"
    "Segmentgruppe_new"       = "Segmentgruppe_FV2504"
    "Segment_new"             = "Segment_FV2504"
    "Datenelement_new"        = "Datenelement_FV2504"
    "Segment ID_new"          = "Segment ID_FV2504"
    "Code_new"                = "Code_FV2504"
    "Qualifier_new"           = "Qualifier_FV2504"
    "Beschreibung_new"        = "Beschreibung_FV2504"
    "Bedingungsausdruck_new"  = "Bedingungsausdruck_FV2504"
    "Bedingung_new"           = "Bedingung_FV2504"
}

for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $current = [string]$cell.Value2
    if ($oldToNew.ContainsKey($current)) {
        $cell.Value2 = $oldToNew[$current]
    }
}

# --- 2) Turn the header + data range into a real worksheet Table. ---------
$tableRange = $ws.Range("A1:U59")
$lo = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $tableRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$lo.Name = "Table1"

# --- 3) Freeze the header row (split below row 1, active pane bottom-left).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
